$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The FARMING (column A) rotation got reshuffled: Hashbrowns/Cabbage/Wheat
# effectively moved down one slot, with Wheat now taking on the "60% - Accent1"
# banding that used to belong to the row below it.
$ws.Range("A2").Value = "Cabbage"

$ws.Range("A3").Value = "Wheat"
$ws.Range("A3").Style = "60% - Accent1"

$ws.Range("A4").Value = "Hashbrowns"

# MAX ENERGY (column J) gains two new cooked goods; the old Chocolate Bar
# entry slides down to make room, and Pancakes takes its old banding.
$ws.Range("J6").Value = "Cranberry Sauce"
$ws.Range("J6").Style = "Accent1"

# New MINING (column D) ingredient used by a recipe in row 5.
$ws.Range("D5").Value = "Stuffing"
$ws.Range("D5").Style = "60% - Accent1"

$ws.Range("J5").Value = "Pancakes"
$ws.Range("J5").Style = "60% - Accent1"

$ws.Range("J7").Value = "Chocolate Bar"
$ws.Range("J7").Style = "Accent1"

# Restore the author's last active selection.
$ws.Range("B6").Select()
